# Modify turn analysis and in-cluster proportion to include multiple phase
# restriction options. This inserts three additional phase-restriction
# columns (C, D, E) after the existing single value column (B), and
# updates column B with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns B (updated), C, D, E (new)
$data = @(
    @{ Row = 1;  B = 12500.0; C = 24500.0; D = 26000.0; E = 32400.0 },
    @{ Row = 2;  B = 15000.0; C = 21000.0; D = 23500.0; E = 32400.0 },
    @{ Row = 3;  B = 8000.0;  C = 16000.0; D = 18000.0; E = 32400.0 },
    @{ Row = 4;  B = 9000.0;  C = 18000.0; D = 19000.0; E = 32400.0 },
    @{ Row = 5;  B = 11600.0; C = 18000.0; D = 20000.0; E = 32400.0 },
    @{ Row = 6;  B = 5000.0;  C = 23500.0; D = 20000.0; E = 32400.0 },
    @{ Row = 7;  B = 12000.0; C = 32400.0; D = 32400.0; E = 32400.0 },
    @{ Row = 8;  B = 10000.0; C = 13500.0; D = 16000.0; E = 32400.0 },
    @{ Row = 9;  B = 12000.0; C = 15800.0; D = 18000.0; E = 32400.0 },
    @{ Row = 10; B = 10000.0; C = 24500.0; D = 27000.0; E = 32400.0 },
    @{ Row = 11; B = 5000.0;  C = 8000.0;  D = 10000.0; E = 32400.0 },
    @{ Row = 12; B = 9000.0;  C = 19000.0; D = 25000.0; E = 32400.0 },
    @{ Row = 13; B = 13000.0; C = 18000.0; D = 22000.0; E = 32400.0 }
)

foreach ($item in $data) {
    $r = $item.Row

    # Copy the existing formatting of column B onto the new C/D/E cells
    # so they share the same style (s="2") as the rest of the numeric
    # column before assigning the new values.
    $ws.Range("B$r").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("E$r").PasteSpecial(-4122)

    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
}
